# Insert 4 new weekly-report rows at the top of the Zanahoria block
# (row 598), pushing the existing rows 598-697 down to 602-701.
# This mirrors the "dimension A1:R697 -> A1:R701" change in the diff:
# four brand-new rows are added, and every previously-existing row
# keeps its original data, just shifted down by four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows before the current row 598.
$ws.Range("598:601").Insert()

# Common/static values shared by every row in this block.
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$categoriaId = 100114013
$categoria   = "Zanahoria"
$variedad    = "Sin especificar"
$unidad      = "`$/saco 20 kilos"
$kgUnidades  = 20
$clasif      = "Hortaliza"

function Set-ZanahoriaRow {
    param(
        $Row,
        $Fecha,
        $Calidad,
        $Volumen,
        $PrecioMin,
        $PrecioMax,
        $PrecioProm,
        $Origen,
        $PrecioKg
    )

    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $categoriaId
    $ws.Cells.Item($Row, 7).Value  = $categoria
    $ws.Cells.Item($Row, 8).Value  = $variedad
    $ws.Cells.Item($Row, 9).Value  = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $kgUnidades
    $ws.Cells.Item($Row, 18).Value = $clasif
}

Set-ZanahoriaRow 598 44505 "Primera" 320 6000 7000 6500 "Chillán"              325
Set-ZanahoriaRow 599 44505 "Primera" 340 6000 7000 6500 "Región Metropolitana" 325
Set-ZanahoriaRow 600 44505 "Segunda" 106 4000 4000 4000 "Chillán"              200
Set-ZanahoriaRow 601 44505 "Segunda" 160 4000 4000 4000 "Región Metropolitana" 200
